# Update country stats and shared-string labels per the upstream CSV refresh
# (countries reordered: Chile/Pakistan block, Irak/Estonia swap, Sierra Leona move;
#  plus same-day case/recovered/death count refreshes, and the footer timestamp).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, optional new country label (column A), and a map of
# column letter -> new numeric value for columns B..H.
$rowUpdates = @(
    @{ Row = 1; Country = 'Datos actualizados a 10 de Abril de 2020 a las 17:52'; Values = @{  } }
    @{ Row = 4; Country = $null; Values = @{ B=475659; C=7093; E=431771; G=1147; H=17838 } }
    @{ Row = 5; Country = $null; Values = @{ B=157053; C=3831; E=85415; G=523; H=15970 } }
    @{ Row = 7; Country = $null; Values = @{ B=119498; C=1263; E=64484 } }
    @{ Row = 11; Country = $null; Values = @{ E=55926; G=1038; H=9016 } }
    @{ Row = 14; Country = $null; Values = @{ B=24469; C=418; E=12882 } }
    @{ Row = 16; Country = $null; Values = @{ B=21243; C=478; E=15423 } }
    @{ Row = 22; Country = $null; Values = @{ E=8939; G=9; H=95 } }
    @{ Row = 24; Country = $null; Values = @{ B=7347; C=622; E=6477 } }
    @{ Row = 26; Country = 'Chile'; Values = @{ B=6501; C=529; D=1274; E=5162; F=360; G=8; H=65 } }
    @{ Row = 29; Country = 'Polonia'; Values = @{ B=5955; C=380; D=318; E=5456; F=160; G=7; H=181 } }
    @{ Row = 30; Country = 'Dinamarca'; Values = @{ B=5819; C=184; D=1773; E=3799; F=113; G=10; H=247 } }
    @{ Row = 31; Country = 'Chequia'; Values = @{ B=5589; C=20; D=309; E=5167 } }
    @{ Row = 32; Country = 'Japon'; Values = @{ B=5530; C=183; D=685; E=4746; F=109; G=0; H=99 } }
    @{ Row = 33; Country = 'Rumania'; Values = @{ B=5467; G=17; H=265 } }
    @{ Row = 34; Country = 'Peru'; Values = @{ B=5256; C=0; D=1438; E=3680; F=124; G=0; H=138 } }
    @{ Row = 35; Country = 'Ecuador'; Values = @{ B=4965; D=339; E=4354; F=139; H=272 } }
    @{ Row = 36; Country = 'Pakistan'; Values = @{ B=4695; C=206; D=727; E=3902; F=45; G=1; H=66 } }
    @{ Row = 64; Country = 'Irak'; Values = @{ B=1279; C=47; D=550; E=659; F=0; G=1; H=70 } }
    @{ Row = 65; Country = 'Estonia'; Values = @{ B=1258; C=51; D=93; E=1141; F=9; H=24 } }
    @{ Row = 84; Country = $null; Values = @{ D=76; E=513 } }
    @{ Row = 86; Country = $null; Values = @{ B=595; C=31; E=532 } }
    @{ Row = 104; Country = $null; Values = @{ B=318; C=4; E=286 } }
    @{ Row = 112; Country = $null; Values = @{ B=234; C=16; E=177 } }
    @{ Row = 140; Country = $null; Values = @{ E=58; G=1; H=3 } }
    @{ Row = 195; Country = 'Sierra Leona'; Values = @{ C=1 } }
    @{ Row = 196; Country = 'Republica de Africa Central'; Values = @{ E=8; H=0 } }
    @{ Row = 197; Country = 'Islas Turcas y Caicos'; Values = @{ D=0; E=7; H=1 } }
    @{ Row = 198; Country = 'Santa Sede'; Values = @{ B=8; D=2; E=6 } }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row
    if ($update.Country -ne $null) {
        $ws.Cells.Item($row, 1).Value = $update.Country
    }
    foreach ($col in $update.Values.Keys) {
        $colIndex = [int][char]$col - [int][char]"A" + 1
        $ws.Cells.Item($row, $colIndex).Value = $update.Values[$col]
    }
}

